$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "57.117.98"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  -1.31%  "

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.984.05"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  -2.39%  "

$ws.Range("E4").Value = "  +0.04%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "500.29"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -4.89%  "

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "137.38"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  -3.80%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -4.26%  "

$ws.Range("E9").Value = "  -5.31%  "

$ws.Range("E10").Value = "  -4.53%  "

$ws.Range("E11").Value = "  -4.21%  "

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.494.18"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  -2.38%  "

$ws.Range("E13").Value = "  -2.57%  "

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "26.06"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  -3.87%  "

$ws.Range("E15").Value = "  -5.47%  "

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "57.158.55"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  -1.19%  "

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.10"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  -2.59%  "

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.982.29"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  -2.46%  "

$ws.Range("E19").Value = "  -3.35%  "

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.86"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -3.11%  "

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "320.40"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  -5.27%  "

$ws.Range("E23").Value = "  +1.03%  "

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.491"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  -2.52%  "

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "63.04"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  -3.01%  "

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  -0.20%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.162"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -5.70%  "

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0893"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  -8.61%  "

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.64"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  -4.21%  "

$ws.Range("E30").Value = "  -3.69%  "

$ws.Range("E31").Value = "  -4.01%  "

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.17"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  -6.09%  "

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "20.15"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  -4.68%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "154.66"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -1.16%  "

$ws.Range("E35").Value = "  -3.59%  "

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.78"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  -3.57%  "

$ws.Range("E37").Value = "  -6.57%  "

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "24.48"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  -6.97%  "

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0663"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  -5.50%  "

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "37.79"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  -0.19%  "

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.015.41"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -2.49%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("E43").Value = "  -3.99%  "

$ws.Range("E44").Value = "  -2.85%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.187.18"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  -5.95%  "

$ws.Range("E46").Value = "  -6.16%  "

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.96"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -1.26%  "

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.929"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -9.93%  "

$ws.Range("E49").Value = "  -4.50%  "

$ws.Range("E50").Value = "  -4.59%  "

$ws.Range("E51").Value = "  -10.92%  "
